# TC03_Canine_Filter_PrimDisSite-Lung.xlsx
# Fixed Diagnosis, FileAssociation, FileFormat, FileType, NeuteredStatus, PrimeDiseaseSite
#
# The "startup" sheet holds, per tab (CasesTab / SamplesTab / FilesTab), the
# Neo4j "query" (column B) used to build that tab. This change drops the
# trailing "Cohort" column from the CasesTab query (column B, row 2) - the
# Cohort node/field is no longer part of that RETURN clause.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$casesTabQuery = @'
MATCH (s:study)<-[*]-(c:case)<--(demo:demographic)
MATCH (c)<--(diag:diagnosis)
OPTIONAL MATCH (samp:sample)-->(c)
OPTIONAL MATCH (co:cohort)<-[*]-(c)
WITH DISTINCT c, s, demo, diag, co
WHERE diag.primary_disease_site IN ['Lung']
RETURN  coalesce(c.case_id, '') AS `Case ID` ,
        coalesce(s.clinical_study_designation, '') AS `Study Code` ,
        coalesce(s.clinical_study_type, '') AS  `Study Type`,
        coalesce(demo.breed, '') AS Breed ,
        coalesce(diag.disease_term, '') AS Diagnosis ,
        coalesce(diag.stage_of_disease, '') AS `Stage of Disease` ,
        coalesce(demo.patient_age_at_enrollment, '') AS Age ,
        coalesce(demo.sex, '') AS Sex ,
        coalesce(demo.neutered_indicator, '') AS `Neutered Status`,
        coalesce(demo.weight, '') AS `Weight (kg)`,
        coalesce(diag.best_response, '') AS `Response to Treatment`
'@

$ws.Range("B2").Value = $casesTabQuery

# Row heights shrink now that the CasesTab query lost its last line.
$ws.Rows.Item(2).RowHeight = 244.8
$ws.Rows.Item(3).RowHeight = 244.8
$ws.Rows.Item(4).RowHeight = 244.8

# Reflect the saved selection/scroll position (top of sheet, B2 selected).
$ws.Activate() | Out-Null
$ws.Range("B2").Select() | Out-Null
